$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column from 45172 to 45175 for all existing data rows (2..527)
$ws.Range("C2:C527").Value = 45175

# Row 527 gains an explicit row height (customHeight) in the saved file
$ws.Range("A527").RowHeight = 15

# Append the new record as row 528
$ws.Range("A528").Value = "A 41189-2023"
$ws.Range("B528").Value = 45174
$ws.Range("C528").Value = 45175
$ws.Range("D528").Value = "VÄRMLANDS LÄN"
$ws.Range("E528").Value = "HAGFORS"
$ws.Range("G528").Value = 3
$ws.Range("H528").Value = 0
$ws.Range("I528").Value = 0
$ws.Range("J528").Value = 0
$ws.Range("K528").Value = 0
$ws.Range("L528").Value = 0
$ws.Range("M528").Value = 0
$ws.Range("N528").Value = 0
$ws.Range("O528").Value = 0
$ws.Range("P528").Value = 0
$ws.Range("Q528").Value = 0
$ws.Range("R528").Value = ""

# Match number formatting used by the rest of column B/C (date style)
$ws.Range("B528:C528").NumberFormat = "YYYY-MM-DD"

# Match the wrap-text style used by the rest of column R
$ws.Range("R528").WrapText = $true
